$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the identifying values (Id, Antal, Ost, Nord) between row 2 and row 4
$cols = @("A", "I", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr4 = "${col}4"
    $val2 = $ws.Range($addr2).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr2).Value2 = $val4
    $ws.Range($addr4).Value2 = $val2
}
